# runScript.xlsx update — restructure the 3-row "test case" table from a
# 9-column (A:I) layout down to a 7-column (A:G) layout, retarget several
# step names/args to a single JSON-ish runScript payload, add a wrapped
# left-aligned style for the new F column, resize columns and row 3, and
# move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the old H:I columns (and their row-3 contents) entirely ----
# Clear() removes both value+style and shrinks the sheet's used range
# (dimension) down to A1:G3, matching the target layout, without leaving
# the old column-width overrides on H/I in a broken state the way a
# column Delete() would.
$ws.Range("H1:I3").Clear() | Out-Null

# --- 2. Row 1 — step names -------------------------------------------
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "open"
$ws.Range("C1").Value = "setWindowSize"
$ws.Range("D1").Value = "assertNotChecked"
$ws.Range("E1").Value = "runScript"
$ws.Range("F1").Value = "runScript"
$ws.Range("G1").Value = "wait"

# --- 3. Row 2 — first argument row ------------------------------------
$ws.Range("A2").ClearContents() | Out-Null
$ws.Range("B2").Value = "https://task.hugang.io/login"
$ws.Range("C2").Value = "945x1012"
$ws.Range("D2").Value = "xpath=//*[@id='username']"
$ws.Range("E2").Value = "xpath=//*[@id='username']"
$ws.Range("F2").ClearContents() | Out-Null
$ws.Range("G2").ClearContents() | Out-Null

# --- 4. Row 3 — second argument / payload row -------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").ClearContents() | Out-Null
$ws.Range("C3").ClearContents() | Out-Null
$ws.Range("D3").ClearContents() | Out-Null
$ws.Range("E3").Value = "arguments[0].value = 'Hello, world!'"
$jsonPayload = "{`n  ""target"": ""xpath=//*[@id='username']"",""value"":""arguments[0].value = 'Hello, world!'""`n}"
$ws.Range("F3").Value = $jsonPayload
$ws.Range("G3").Value = 1000

# --- 5. New style for F2/F3: left+center aligned, wrapped, bordered ----
$ws.Range("F2:F3").HorizontalAlignment = -4131   # xlLeft
$ws.Range("F2:F3").VerticalAlignment = -4108     # xlCenter
$ws.Range("F2:F3").WrapText = $true

# --- 6. Column widths: F grows to fit the JSON text, G shrinks --------
$ws.Columns.Item(6).ColumnWidth = 38
$ws.Columns.Item(7).ColumnWidth = 4.2857142857142856

# --- 7. Row 3 grows tall enough to show the wrapped JSON text ---------
$ws.Rows.Item(3).RowHeight = 78.75

# --- 8. Move the saved selection cursor --------------------------------
$ws.Range("E8").Select() | Out-Null
